$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# Update Version and Date values
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row 11 ("Jurisdiction") after "Contact" (row 10) and before "Description" (old row 11)
$meta.Rows.Item(11).Insert()
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# --- Sheet "Elements" ---
$els = $wb.Worksheets.Item("Elements")

# Add the II-1 constraint text to ExternalProcedure.typeId (row 5), column AJ (Constraint(s))
$els.Range("AJ5").Value = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}`n"

# Add a new row 17 for ExternalProcedure.sdtcAuthor, matching formatting of row 16
$els.Range("A16:AK16").Copy()
$els.Range("A17:AK17").PasteSpecial(-4122)

$els.Range("A17").Value = "ExternalProcedure.sdtcAuthor"
$els.Range("B17").Value = "ExternalProcedure.sdtcAuthor"
$els.Range("F17").Value = "0"
$els.Range("G17").Value = "*"
$els.Range("K17").Value = "http://hl7.org/cda/stds/core/StructureDefinition/Author`n"
$els.Range("AF17").Value = "ExternalProcedure.sdtcAuthor"
$els.Range("AG17").Value = "0"
$els.Range("AH17").Value = "*"
